$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 2).Value = 0.5089696350451618
$ws.Cells.Item(2, 3).Value = 0.3725824631726109
$ws.Cells.Item(2, 5).Value = 0.2833489666963693
$ws.Cells.Item(2, 6).Value = 2.068094517207555
$ws.Cells.Item(2, 7).Value = 0.7770739020213142
$ws.Cells.Item(2, 8).Value = 0.8777017147733304
$ws.Cells.Item(2, 10).Value = 0.07059028345345553
$ws.Cells.Item(2, 11).Value = 0.2334338854878126
$ws.Cells.Item(2, 12).Value = 0.4213504119095859
$ws.Cells.Item(2, 13).Value = 0.215225388113204
$ws.Cells.Item(2, 15).Value = 3.323135395211011
$ws.Cells.Item(3, 2).Value = 0.4715632714483036
$ws.Cells.Item(3, 3).Value = 0.3750206972920722
$ws.Cells.Item(3, 5).Value = 0.2835279461376246
$ws.Cells.Item(3, 6).Value = 2.073785284769158
$ws.Cells.Item(3, 7).Value = 0.7839489105253818
$ws.Cells.Item(3, 8).Value = 0.8845694471904721
$ws.Cells.Item(3, 10).Value = 0.06840293147484644
$ws.Cells.Item(3, 11).Value = 0.2036824312599208
$ws.Cells.Item(3, 12).Value = 0.4157219087694983
$ws.Cells.Item(3, 13).Value = 0.2066109941149605
$ws.Cells.Item(3, 15).Value = 3.35184895133284
$ws.Cells.Item(4, 2).Value = 0.4486837299064916
$ws.Cells.Item(4, 3).Value = 0.3766016061303468
$ws.Cells.Item(4, 5).Value = 0.2837345105077915
$ws.Cells.Item(4, 6).Value = 2.078234508088492
$ws.Cells.Item(4, 7).Value = 0.7885959409305201
$ws.Cells.Item(4, 8).Value = 0.8891051468936695
$ws.Cells.Item(4, 10).Value = 0.06704894626164659
$ws.Cells.Item(4, 11).Value = 0.1853510153708555
$ws.Cells.Item(4, 12).Value = 0.4124128509810987
$ws.Cells.Item(4, 13).Value = 0.2013867044975868
$ws.Cells.Item(4, 15).Value = 3.371041113095885
$ws.Cells.Item(5, 2).Value = 0.4393830043618436
$ws.Cells.Item(5, 3).Value = 0.3772669542948286
$ws.Cells.Item(5, 5).Value = 0.2838430591478591
$ws.Cells.Item(5, 6).Value = 2.080288073649221
$ws.Cells.Item(5, 7).Value = 0.790596705563523
$ws.Cells.Item(5, 8).Value = 0.8910337465252738
$ws.Cells.Item(5, 10).Value = 0.06649446100686163
$ws.Cells.Item(5, 11).Value = 0.1778652393571605
$ws.Cells.Item(5, 12).Value = 0.4111014987260972
$ws.Cells.Item(5, 13).Value = 0.1992742932564937
$ws.Cells.Item(5, 15).Value = 3.379255002289199
$ws.Cells.Item(6, 2).Value = 0.4378400285573036
$ws.Cells.Item(6, 3).Value = 0.3773787112322413
$ws.Cells.Item(6, 5).Value = 0.2838625574397646
$ws.Cells.Item(6, 6).Value = 2.080643599635529
$ws.Cells.Item(6, 7).Value = 0.7909353974172362
$ws.Cells.Item(6, 8).Value = 0.8913588391728595
$ws.Cells.Item(6, 10).Value = 0.06640222532357853
$ws.Cells.Item(6, 11).Value = 0.1766213058097179
$ws.Cells.Item(6, 12).Value = 0.4108859968205394
$ws.Cells.Item(6, 13).Value = 0.1989245328571236
$ws.Cells.Item(6, 15).Value = 3.380642649168948
$ws.Cells.Item(7, 2).Value = 0.4485582035053142
$ws.Cells.Item(7, 3).Value = 0.3766104937062984
$ws.Cells.Item(7, 5).Value = 0.2837358756669595
$ws.Cells.Item(7, 6).Value = 2.078261229121019
$ws.Cells.Item(7, 7).Value = 0.788622490446194
$ws.Cells.Item(7, 8).Value = 0.8891308315666819
$ws.Cells.Item(7, 10).Value = 0.06704147927942472
$ws.Cells.Item(7, 11).Value = 0.1852501219855895
$ws.Cells.Item(7, 12).Value = 0.4123950151003442
$ws.Cells.Item(7, 13).Value = 0.2013581486042462
$ws.Cells.Item(7, 15).Value = 3.371150297305363
$ws.Cells.Item(8, 2).Value = 0.4960540490984044
$ws.Cells.Item(8, 3).Value = 0.3734057911422308
$ws.Cells.Item(8, 5).Value = 0.2833906524924465
$ws.Cells.Item(8, 6).Value = 2.069858652818596
$ws.Cells.Item(8, 7).Value = 0.779356043642295
$ws.Cells.Item(8, 8).Value = 0.8800035722449948
$ws.Cells.Item(8, 10).Value = 0.06983837070048793
$ws.Cells.Item(8, 11).Value = 0.2231891338962697
$ws.Cells.Item(8, 12).Value = 0.4193793305002345
$ws.Cells.Item(8, 13).Value = 0.2122417679332287
$ws.Cells.Item(8, 15).Value = 3.332711754468491
$ws.Cells.Item(9, 2).Value = 0.5898657433008623
$ws.Cells.Item(9, 3).Value = 0.3677847633633622
$ws.Cells.Item(9, 5).Value = 0.2834781154435078
$ws.Cells.Item(9, 6).Value = 2.060947110676466
$ws.Cells.Item(9, 7).Value = 0.7645626694107079
$ws.Cells.Item(9, 8).Value = 0.8646316951452846
$ws.Cells.Item(9, 10).Value = 0.07523532708048464
$ws.Cells.Item(9, 11).Value = 0.2970621964340125
$ws.Cells.Item(9, 12).Value = 0.4342337467224553
$ws.Cells.Item(9, 13).Value = 0.2340930269317241
$ws.Cells.Item(9, 15).Value = 3.269719825073366
$ws.Cells.Item(10, 2).Value = 0.6591703046944417
$ws.Cells.Item(10, 3).Value = 0.3640569675008276
$ws.Cells.Item(10, 5).Value = 0.2840052083290274
$ws.Cells.Item(10, 6).Value = 2.058998047223355
$ws.Cells.Item(10, 7).Value = 0.7557536034308256
$ws.Cells.Item(10, 8).Value = 0.8548734037537571
$ws.Cells.Item(10, 10).Value = 0.07914609643337656
$ws.Cells.Item(10, 11).Value = 0.3509972601552249
$ws.Cells.Item(10, 12).Value = 0.44584526976017
$ws.Cells.Item(10, 13).Value = 0.2504493659439362
$ws.Cells.Item(10, 15).Value = 3.23098109710763
$ws.Cells.Item(11, 2).Value = 0.6907758172211231
$ws.Cells.Item(11, 3).Value = 0.3624478619589251
$ws.Cells.Item(11, 5).Value = 0.2843448389981944
$ws.Cells.Item(11, 6).Value = 2.05910698439466
$ws.Cells.Item(11, 7).Value = 0.7521934378389616
$ws.Cells.Item(11, 8).Value = 0.8507664796376488
$ws.Cells.Item(11, 10).Value = 0.08091323077661627
$ws.Cells.Item(11, 11).Value = 0.3754562422358561
$ws.Cells.Item(11, 12).Value = 0.4512775891581668
$ws.Cells.Item(11, 13).Value = 0.2579543688902888
$ws.Cells.Item(11, 15).Value = 3.214993488533523
$ws.Cells.Item(12, 2).Value = 0.7027546611865318
$ws.Cells.Item(12, 3).Value = 0.3618509636337777
$ws.Cells.Item(12, 5).Value = 0.2844877472274838
$ws.Cells.Item(12, 6).Value = 2.059291122776472
$ws.Cells.Item(12, 7).Value = 0.7509095982719529
$ws.Cells.Item(12, 8).Value = 0.8492589840950586
$ws.Cells.Item(12, 10).Value = 0.08158066552116594
$ws.Cells.Item(12, 11).Value = 0.3847068125993189
$ws.Cells.Item(12, 12).Value = 0.4533560949781759
$ws.Cells.Item(12, 13).Value = 0.2608054138954756
$ws.Cells.Item(12, 15).Value = 3.20917434493191
$ws.Cells.Item(13, 2).Value = 0.7001743470515862
$ws.Cells.Item(13, 3).Value = 0.3619789638001656
$ws.Cells.Item(13, 5).Value = 0.2844563344763635
$ws.Cells.Item(13, 6).Value = 2.059245114790897
$ws.Cells.Item(13, 7).Value = 0.7511832351909433
$ws.Cells.Item(13, 8).Value = 0.8495815293549285
$ws.Cells.Item(13, 10).Value = 0.08143699940758609
$ws.Cells.Item(13, 11).Value = 0.3827150574276459
$ws.Cells.Item(13, 12).Value = 0.4529075026413807
$ws.Cells.Item(13, 13).Value = 0.2601909904421262
$ws.Cells.Item(13, 15).Value = 3.210417150310946
$ws.Cells.Item(14, 2).Value = 0.6917611174391141
$ws.Cells.Item(14, 3).Value = 0.3623985057326777
$ws.Cells.Item(14, 5).Value = 0.2843563099529369
$ws.Cells.Item(14, 6).Value = 2.059119271504116
$ws.Cells.Item(14, 7).Value = 0.7520865264334446
$ws.Cells.Item(14, 8).Value = 0.8506415012762147
$ws.Cells.Item(14, 10).Value = 0.08096817609426665
$ws.Cells.Item(14, 11).Value = 0.3762175260781362
$ws.Cells.Item(14, 12).Value = 0.4514481612031034
$ws.Cells.Item(14, 13).Value = 0.2581887455854925
$ws.Cells.Item(14, 15).Value = 3.214510034403943
$ws.Cells.Item(15, 2).Value = 0.6866091188982182
$ws.Cells.Item(15, 3).Value = 0.3626571057856136
$ws.Cells.Item(15, 5).Value = 0.2842969021760915
$ws.Cells.Item(15, 6).Value = 2.059060788641716
$ws.Cells.Item(15, 7).Value = 0.7526481951378301
$ws.Cells.Item(15, 8).Value = 0.8512969758809135
$ws.Cells.Item(15, 10).Value = 0.08068078068964724
$ws.Cells.Item(15, 11).Value = 0.3722360849874633
$ws.Cells.Item(15, 12).Value = 0.4505570542384305
$ws.Cells.Item(15, 13).Value = 0.2569634871528947
$ws.Cells.Item(15, 15).Value = 3.217047648125202
$ws.Cells.Item(16, 2).Value = 0.6571063243906394
$ws.Cells.Item(16, 3).Value = 0.3641638680884629
$ws.Cells.Item(16, 5).Value = 0.2839850167289164
$ws.Cells.Item(16, 6).Value = 2.059010937857423
$ws.Cells.Item(16, 7).Value = 0.7559952677836819
$ws.Cells.Item(16, 8).Value = 0.8551484800454148
$ws.Cells.Item(16, 10).Value = 0.07903036811700304
$ws.Cells.Item(16, 11).Value = 0.3493972259228713
$ws.Cells.Item(16, 12).Value = 0.4454932605561197
$ws.Cells.Item(16, 13).Value = 0.2499601755417515
$ws.Cells.Item(16, 15).Value = 3.2320588159522
$ws.Cells.Item(17, 2).Value = 0.6390268789011202
$ws.Cells.Item(17, 3).Value = 0.3651103988632762
$ws.Cells.Item(17, 5).Value = 0.2838192163259059
$ws.Cells.Item(17, 6).Value = 2.059235150612054
$ws.Cells.Item(17, 7).Value = 0.758163123223774
$ws.Cells.Item(17, 8).Value = 0.8575962885482511
$ws.Cells.Item(17, 10).Value = 0.0780148255222386
$ws.Cells.Item(17, 11).Value = 0.3353663886092306
$ws.Cells.Item(17, 12).Value = 0.4424251156054311
$ws.Cells.Item(17, 13).Value = 0.2456802280182018
$ws.Cells.Item(17, 15).Value = 3.241686351290724
$ws.Cells.Item(18, 2).Value = 0.6286354965614009
$ws.Cells.Item(18, 3).Value = 0.3656629795688584
$ws.Cells.Item(18, 5).Value = 0.2837332546671512
$ws.Cells.Item(18, 6).Value = 2.05945783676961
$ws.Cells.Item(18, 7).Value = 0.7594521003439851
$ws.Cells.Item(18, 8).Value = 0.8590354722821658
$ws.Cells.Item(18, 10).Value = 0.07742959455433862
$ws.Cells.Item(18, 11).Value = 0.3272890677203009
$ws.Cells.Item(18, 12).Value = 0.4406745468027111
$ws.Cells.Item(18, 13).Value = 0.2432245911110513
$ws.Cells.Item(18, 15).Value = 3.247377737413217
$ws.Cells.Item(19, 2).Value = 0.6251184537301526
$ws.Cells.Item(19, 3).Value = 0.3658514765388166
$ws.Cells.Item(19, 5).Value = 0.2837057662581906
$ws.Cells.Item(19, 6).Value = 2.059549340223384
$ws.Cells.Item(19, 7).Value = 0.7598957531357726
$ws.Cells.Item(19, 8).Value = 0.8595281274374216
$ws.Cells.Item(19, 10).Value = 0.07723125435560974
$ws.Cells.Item(19, 11).Value = 0.3245530132597594
$ws.Cells.Item(19, 12).Value = 0.4400842698768201
$ws.Cells.Item(19, 13).Value = 0.24239420427174
$ws.Cells.Item(19, 15).Value = 3.249331174558847
$ws.Cells.Item(20, 2).Value = 0.6409507026596089
$ws.Cells.Item(20, 3).Value = 0.3650087946105369
$ws.Cells.Item(20, 5).Value = 0.2838358933986953
$ws.Cells.Item(20, 6).Value = 2.059201584603613
$ws.Cells.Item(20, 7).Value = 0.757927995611773
$ws.Cells.Item(20, 8).Value = 0.857332479351065
$ws.Cells.Item(20, 10).Value = 0.07812304766473233
$ws.Cells.Item(20, 11).Value = 0.3368607395351262
$ws.Cells.Item(20, 12).Value = 0.4427502620950463
$ws.Cells.Item(20, 13).Value = 0.246135208379286
$ws.Cells.Item(20, 15).Value = 3.240645557970012
$ws.Cells.Item(21, 2).Value = 0.6942320074623467
$ws.Cells.Item(21, 3).Value = 0.3622749388547195
$ws.Cells.Item(21, 5).Value = 0.2843853020397482
$ws.Cells.Item(21, 6).Value = 2.05915235904483
$ws.Cells.Item(21, 7).Value = 0.7518194621566963
$ws.Cells.Item(21, 8).Value = 0.8503288675245102
$ws.Cells.Item(21, 10).Value = 0.08110592834913177
$ws.Cells.Item(21, 11).Value = 0.3781263253268321
$ws.Cells.Item(21, 12).Value = 0.4518762257746545
$ws.Cells.Item(21, 13).Value = 0.2587766090580388
$ws.Cells.Item(21, 15).Value = 3.213301477405508
$ws.Cells.Item(22, 2).Value = 0.729115396584092
$ws.Cells.Item(22, 3).Value = 0.3605606729887398
$ws.Cells.Item(22, 5).Value = 0.2848276804667584
$ws.Cells.Item(22, 6).Value = 2.059952908581408
$ws.Cells.Item(22, 7).Value = 0.7482020861764411
$ws.Cells.Item(22, 8).Value = 0.8460296512664414
$ws.Cells.Item(22, 10).Value = 0.08304525318536093
$ws.Cells.Item(22, 11).Value = 0.4050284032836942
$ws.Cells.Item(22, 12).Value = 0.4579652428279815
$ws.Cells.Item(22, 13).Value = 0.2670912299063701
$ws.Cells.Item(22, 15).Value = 3.196800364979211
$ws.Cells.Item(23, 2).Value = 0.7104921495756003
$ws.Cells.Item(23, 3).Value = 0.3614689882117794
$ws.Cells.Item(23, 5).Value = 0.2845839723385986
$ws.Cells.Item(23, 6).Value = 2.059449537784303
$ws.Cells.Item(23, 7).Value = 0.7500984368964723
$ws.Cells.Item(23, 8).Value = 0.8482988015877311
$ws.Cells.Item(23, 10).Value = 0.08201113870555332
$ws.Cells.Item(23, 11).Value = 0.3906765989455039
$ws.Cells.Item(23, 12).Value = 0.454704075504921
$ws.Cells.Item(23, 13).Value = 0.2626488012672468
$ws.Cells.Item(23, 15).Value = 3.205482003428955
$ws.Cells.Item(24, 2).Value = 0.6400809333715358
$ws.Cells.Item(24, 3).Value = 0.3650547037024232
$ws.Cells.Item(24, 5).Value = 0.2838283245367847
$ws.Cells.Item(24, 6).Value = 2.059216467649321
$ws.Cells.Item(24, 7).Value = 0.7580341639611774
$ws.Cells.Item(24, 8).Value = 0.8574516480989161
$ws.Cells.Item(24, 10).Value = 0.07807412473643183
$ws.Cells.Item(24, 11).Value = 0.3361851770996793
$ws.Cells.Item(24, 12).Value = 0.4426032217738936
$ws.Cells.Item(24, 13).Value = 0.2459294962873813
$ws.Cells.Item(24, 15).Value = 3.241115613487509
$ws.Cells.Item(25, 2).Value = 0.5644181898633747
$ws.Cells.Item(25, 3).Value = 0.3692346344538358
$ws.Cells.Item(25, 5).Value = 0.2833729454683436
$ws.Cells.Item(25, 6).Value = 2.062549395341037
$ws.Cells.Item(25, 7).Value = 0.7682029936467316
$ws.Cells.Item(25, 8).Value = 0.8685202067709668
$ws.Cells.Item(25, 10).Value = 0.07378479184326636
$ws.Cells.Item(25, 11).Value = 0.2771358052026187
$ws.Cells.Item(25, 12).Value = 0.4300920765904408
$ws.Cells.Item(25, 13).Value = 0.2281280481998991
$ws.Cells.Item(25, 15).Value = 3.285435779323421